# Daily attendance processing - 2025-10-18 06:25:02
# Swap the order of names/emails listed in the "Recorded By" column (G)
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows whose "Recorded By" (column G) value is "System, dnasr281@gmail.com"
# and needs to become "dnasr281@gmail.com, System"
$systemRows = @(3,6,11,12,13,14,15,17,30,33,38,39,40,41,42,44,57,60,65,66,67,68,69,71,87,88,89,93,95,96,99,113,114,115,119,121,122,125,139,140,141,145,147,148,151)

foreach ($r in $systemRows) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value() -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}

# Rows whose "Recorded By" (column G) value is "admin@admin.com, dnasr281@gmail.com"
# and needs to become "dnasr281@gmail.com, admin@admin.com"
$adminRows = @(90,116,142)

foreach ($r in $adminRows) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value() -eq "admin@admin.com, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, admin@admin.com"
    }
}
